$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 14504.75
$ws.Range("I18").Value = 10762.5
$ws.Range("J18").Value = 16999.584
$ws.Range("K18").Value = 10762.5
$ws.Range("L18").Value = 16999.584
$ws.Range("M18").Value = -10478.5
$ws.Range("N18").Value = -17567.584
$ws.Range("H62").Value = 9195.286
$ws.Range("J62").Value = 10659.667
$ws.Range("L62").Value = 10659.667
$ws.Range("N62").Value = -11907.667
$ws.Range("H65").Value = 9195.286
$ws.Range("J65").Value = 10659.667
$ws.Range("L65").Value = 53298.335
$ws.Range("N65").Value = -59538.335
$ws.Range("H98").Value = 3227.8235
$ws.Range("I98").Value = 2749.5
$ws.Range("K98").Value = 2749.5
$ws.Range("M98").Value = -1251.5
$ws.Range("H100").Value = 2178.2856
$ws.Range("I100").Value = 1452.5
$ws.Range("J100").Value = 3146
$ws.Range("K100").Value = 1452.5
$ws.Range("L100").Value = 3146
$ws.Range("M100").Value = -911.5
$ws.Range("N100").Value = -4228
$ws.Range("H107").Value = 588.26666
$ws.Range("I107").Value = 495.6
$ws.Range("K107").Value = 495.6
$ws.Range("M107").Value = 1424.4
$ws.Range("H116").Value = 7500
$ws.Range("J116").Value = 7500
$ws.Range("L116").Value = 7500
$ws.Range("N116").Value = -14384
$ws.Range("H121").Value = 1047.25
$ws.Range("J121").Value = 1051.5454
$ws.Range("L121").Value = 3154.6362
$ws.Range("N121").Value = -6648.6362
$ws.Range("H122").Value = 3227.8235
$ws.Range("I122").Value = 2749.5
$ws.Range("K122").Value = 8248.5
$ws.Range("M122").Value = -5798.5
$ws.Range("H125").Value = 1080
$ws.Range("I125").Value = 1080
$ws.Range("K125").Value = 9720
$ws.Range("M125").Value = -7260
$ws.Range("H129").Value = 980.3103599999999
$ws.Range("J129").Value = 989.95
$ws.Range("L129").Value = 2969.85
$ws.Range("N129").Value = -12969.85
$ws.Range("H132").Value = 817.74286
$ws.Range("I132").Value = 806.69696
$ws.Range("K132").Value = 2420.09088
$ws.Range("M132").Value = 109.9091200000003
$ws.Range("H135").Value = 581.7778
$ws.Range("I135").Value = 462.42856
$ws.Range("J135").Value = 999.5
$ws.Range("K135").Value = 4161.85704
$ws.Range("L135").Value = 8995.5
$ws.Range("M135").Value = -1626.85704
$ws.Range("N135").Value = -14065.5
$ws.Range("H137").Value = 2589.8262
$ws.Range("I137").Value = 2134
$ws.Range("J137").Value = 2658.2
$ws.Range("K137").Value = 6402
$ws.Range("L137").Value = 7974.599999999999
$ws.Range("M137").Value = -3852
$ws.Range("N137").Value = -13074.6
$ws.Range("H138").Value = 3457.9546
$ws.Range("I138").Value = 4697.727
$ws.Range("J138").Value = 2218.182
$ws.Range("K138").Value = 14093.181
$ws.Range("L138").Value = 6654.545999999999
$ws.Range("M138").Value = -8953.181
$ws.Range("N138").Value = -16934.546
$ws.Range("H141").Value = 1169327
$ws.Range("I141").Value = 1557880.5
$ws.Range("J141").Value = 3666.3333
$ws.Range("K141").Value = 4673641.5
$ws.Range("L141").Value = 10998.9999
$ws.Range("M141").Value = -4668461.5
$ws.Range("N141").Value = -21358.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3004.8225
$ws.Range("I32").Value = 2308.6458
$ws.Range("K32").Value = 2308.6458
$ws.Range("M32").Value = -2021.6458
$ws.Range("H45").Value = 2320.6365
$ws.Range("J45").Value = 2754.85
$ws.Range("L45").Value = 2754.85
$ws.Range("N45").Value = -3508.85
$ws.Range("H46").Value = 2000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2638
$ws.Range("H61").Value = 4285.4287
$ws.Range("I61").Value = 1999.3334
$ws.Range("K61").Value = 1999.3334
$ws.Range("M61").Value = -1787.3334
$ws.Range("H119").Value = 57498
$ws.Range("J119").Value = 57498
$ws.Range("L119").Value = 57498
$ws.Range("N119").Value = -67174
$ws.Range("H122").Value = 1048.3334
$ws.Range("I122").Value = 972.625
$ws.Range("J122").Value = 1199.75
$ws.Range("K122").Value = 2917.875
$ws.Range("L122").Value = 3599.25
$ws.Range("M122").Value = -467.875
$ws.Range("N122").Value = -8499.25
$ws.Range("H123").Value = 64499
$ws.Range("J123").Value = 64499
$ws.Range("L123").Value = 64499
$ws.Range("N123").Value = -74299
$ws.Range("H136").Value = 4285.4287
$ws.Range("I136").Value = 1999.3334
$ws.Range("K136").Value = 5998.0002
$ws.Range("M136").Value = -3448.0002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 119151.7
$ws.Range("I86").Value = 1438.6
$ws.Range("K86").Value = 1438.6
$ws.Range("M86").Value = -315.5999999999999
$ws.Range("H89").Value = 119151.7
$ws.Range("I89").Value = 1438.6
$ws.Range("K89").Value = 7193
$ws.Range("M89").Value = -1577
$ws.Range("H107").Value = 2930
$ws.Range("I107").Value = 2934.6924
$ws.Range("K107").Value = 2934.6924
$ws.Range("M107").Value = -1014.6924
$ws.Range("H122").Value = 45780
$ws.Range("J122").Value = 45780
$ws.Range("L122").Value = 45780
$ws.Range("N122").Value = -55580
$ws.Range("H134").Value = 295
$ws.Range("I134").Value = 295
$ws.Range("K134").Value = 885
$ws.Range("M134").Value = 1650

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1257.7021
$ws.Range("I31").Value = 915.8889
$ws.Range("J31").Value = 1469.862
$ws.Range("K31").Value = 915.8889
$ws.Range("L31").Value = 1469.862
$ws.Range("M31").Value = -620.8889
$ws.Range("N31").Value = -2059.862
$ws.Range("H34").Value = 1257.7021
$ws.Range("I34").Value = 915.8889
$ws.Range("J34").Value = 1469.862
$ws.Range("K34").Value = 915.8889
$ws.Range("L34").Value = 1469.862
$ws.Range("M34").Value = -713.8889
$ws.Range("N34").Value = -1873.862
$ws.Range("H41").Value = 17033.572
$ws.Range("I41").Value = 12247
$ws.Range("J41").Value = 29000
$ws.Range("K41").Value = 12247
$ws.Range("L41").Value = 29000
$ws.Range("M41").Value = -11819
$ws.Range("N41").Value = -29856
$ws.Range("H50").Value = 17980
$ws.Range("J50").Value = 17980
$ws.Range("L50").Value = 17980
$ws.Range("N50").Value = -19230
$ws.Range("H51").Value = 33000
$ws.Range("J51").Value = 33000
$ws.Range("L51").Value = 33000
$ws.Range("N51").Value = -34472
$ws.Range("H59").Value = 38966.332
$ws.Range("J59").Value = 38966.332
$ws.Range("L59").Value = 38966.332
$ws.Range("N59").Value = -41256.332
$ws.Range("H61").Value = 33000
$ws.Range("J61").Value = 33000
$ws.Range("L61").Value = 33000
$ws.Range("N61").Value = -33696
$ws.Range("H132").Value = 3628.0715
$ws.Range("I132").Value = 3079.6
$ws.Range("K132").Value = 9238.799999999999
$ws.Range("M132").Value = -6708.799999999999
$ws.Range("H141").Value = 69829.25
$ws.Range("J141").Value = 69829.25
$ws.Range("L141").Value = 69829.25
$ws.Range("N141").Value = -80189.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1871.6727
$ws.Range("I68").Value = 799.6667
$ws.Range("J68").Value = 1933.5193
$ws.Range("K68").Value = 2399.0001
$ws.Range("L68").Value = 5800.5579
$ws.Range("M68").Value = -1588.0001
$ws.Range("N68").Value = -7422.5579
$ws.Range("H71").Value = 1871.6727
$ws.Range("I71").Value = 799.6667
$ws.Range("J71").Value = 1933.5193
$ws.Range("K71").Value = 7197.0003
$ws.Range("L71").Value = 17401.6737
$ws.Range("M71").Value = -3141.0003
$ws.Range("N71").Value = -25513.6737
$ws.Range("H103").Value = 18282.2
$ws.Range("I103").Value = 25961.8
$ws.Range("J103").Value = 2923
$ws.Range("K103").Value = 77885.39999999999
$ws.Range("L103").Value = 8769
$ws.Range("M103").Value = -77006.39999999999
$ws.Range("N103").Value = -10527
$ws.Range("H131").Value = 10654378
$ws.Range("J131").Value = 20317.73
$ws.Range("L131").Value = 60953.19
$ws.Range("N131").Value = -71033.19
$ws.Range("H139").Value = 14446.25
$ws.Range("I139").Value = 14446.25
$ws.Range("K139").Value = 43338.75
$ws.Range("M139").Value = -38198.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6300
$ws.Range("J80").Value = 1900
$ws.Range("L80").Value = 1900
$ws.Range("N80").Value = -3896
$ws.Range("H83").Value = 6300
$ws.Range("J83").Value = 1900
$ws.Range("L83").Value = 9500
$ws.Range("N83").Value = -19484
$ws.Range("H102").Value = 2869.9
$ws.Range("I102").Value = 2967.2222
$ws.Range("K102").Value = 2967.2222
$ws.Range("M102").Value = -1345.2222
$ws.Range("H122").Value = 2361.818
$ws.Range("I122").Value = 1157.1428
$ws.Range("K122").Value = 3471.4284
$ws.Range("M122").Value = -1021.4284
$ws.Range("H126").Value = 1918491.1
$ws.Range("I126").Value = 2648209.5
$ws.Range("J126").Value = 2980.125
$ws.Range("K126").Value = 7944628.5
$ws.Range("L126").Value = 8940.375
$ws.Range("M126").Value = -7942158.5
$ws.Range("N126").Value = -13880.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4334
$ws.Range("I7").Value = 2867.1667
$ws.Range("J7").Value = 5311.8887
$ws.Range("K7").Value = 2867.1667
$ws.Range("L7").Value = 5311.8887
$ws.Range("M7").Value = -2755.1667
$ws.Range("N7").Value = -5535.8887
$ws.Range("H16").Value = 8122
$ws.Range("I16").Value = 8785.556
$ws.Range("J16").Value = 2150
$ws.Range("K16").Value = 8785.556
$ws.Range("L16").Value = 2150
$ws.Range("M16").Value = -8615.556
$ws.Range("N16").Value = -2490
$ws.Range("H40").Value = 6829.4
$ws.Range("I40").Value = 3327.8572
$ws.Range("K40").Value = 3327.8572
$ws.Range("M40").Value = -3191.8572
$ws.Range("H100").Value = 2662.3635
$ws.Range("I100").Value = 2525.75
$ws.Range("K100").Value = 2525.75
$ws.Range("M100").Value = -1984.75
$ws.Range("H108").Value = 65000
$ws.Range("J108").Value = 65000
$ws.Range("L108").Value = 65000
$ws.Range("N108").Value = -72680
$ws.Range("H122").Value = 7206
$ws.Range("J122").Value = 14998.5
$ws.Range("L122").Value = 44995.5
$ws.Range("N122").Value = -49895.5
$ws.Range("H126").Value = 4334
$ws.Range("I126").Value = 2867.1667
$ws.Range("J126").Value = 5311.8887
$ws.Range("K126").Value = 8601.500100000001
$ws.Range("L126").Value = 15935.6661
$ws.Range("M126").Value = -6131.500100000001
$ws.Range("N126").Value = -20875.6661
$ws.Range("H132").Value = 3839.818
$ws.Range("J132").Value = 4073.8
$ws.Range("L132").Value = 12221.4
$ws.Range("N132").Value = -17281.4
$ws.Range("H136").Value = 5937.4614
$ws.Range("I136").Value = 3879.4
$ws.Range("K136").Value = 11638.2
$ws.Range("M136").Value = -9088.200000000001
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 5200
$ws.Range("I32").Value = 5200
$ws.Range("K32").Value = 5200
$ws.Range("M32").Value = -4883
$ws.Range("H62").Value = 2066.6667
$ws.Range("I62").Value = 2100
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 2100
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -1476
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 2066.6667
$ws.Range("I65").Value = 2100
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 10500
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -7380
$ws.Range("N65").Value = -16240
$ws.Range("H122").Value = 79755.836
$ws.Range("I122").Value = 118777.44
$ws.Range("J122").Value = 1712.625
$ws.Range("K122").Value = 356332.32
$ws.Range("L122").Value = 5137.875
$ws.Range("M122").Value = -353882.32
$ws.Range("N122").Value = -10037.875
$ws.Range("H132").Value = 4999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 55560604
$ws.Range("I136").Value = 138895890
$ws.Range("J136").Value = 3741.6667
$ws.Range("K136").Value = 416687670
$ws.Range("L136").Value = 11225.0001
$ws.Range("M136").Value = -416685120
$ws.Range("N136").Value = -16325.0001
